$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The template previously had a spare "id" / "{{index}}" mapping row (row 3)
# that duplicated the purpose of the other rows. Remove it so the real
# mapping rows shift up into place.
$ws.Rows("3:3").Delete()

# After the delete, the rows that used to be 4..13 are now 3..12.
# Column E in this template mirrors column B's placeholder text (it used to
# be done with formulas like "=B4" / "=B5+0", which is why several of them
# evaluated to #VALUE!). Replace those formulas with the same literal
# placeholder text that's already in column B, completing the mapping.
$ws.Range("E3").Value = "{{fullName}}"
$ws.Range("E4").Value = "{{oldIndex}}"
$ws.Range("E5").Value = "{{newIndex}}"
$ws.Range("E7").Value = "{{unitPrice}}"

# Rows 6 and 9 (the old B7/B10 "Điện tiêu thụ" / "Tổng tiền thanh toán"
# computed rows) become placeholders too, so the template can be filled in
# directly instead of trying to compute from blank inputs.
$ws.Range("B6").Value = "{{unitsInMonth}}"
$ws.Range("E6").Value = "{{unitsInMonth}}"
$ws.Range("B9").Value = "{{totalPayment}}"
$ws.Range("E9").Value = "{{totalPayment}}"

# Row 8 (Công ghi điện) keeps its B8 -> E8 formula/value link untouched.

# Move the selection, matching the saved workbook state.
$ws.Range("M1").Select()

# The vertical divider line that used to span down to the (now removed)
# last blank row shrinks back up by one row so its anchor still lands on
# the same row it always did relative to the bottom of the sheet.
$shp = $ws.Shapes.Item(1)
$shp.Height = 378.55
